# 🔄 MAJ automatique BRVM via GitHub Actions
#
# Refresh the "Recommandations" and "Top_YTD" sheets with the latest
# automated BRVM computation: updated hausse/baisse day counts, variation
# totals/last-variation percentages, a handful of re-ranked / renamed
# rows, three brand-new rows appended at the bottom of Recommandations,
# and refreshed YTD progression figures (with one ticker dropping out of
# the Top_YTD ranking).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Recommandations
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Recommandations")

# Row 2 - BRVM - SERVICES PUBLICS
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 8
$ws.Cells.Item(2,4).Value = 3223.43
$ws.Cells.Item(2,5).Value = 100.18

# Row 3 - CFAO MOTORS CI
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = 4
$ws.Cells.Item(3,4).Value = 2595
$ws.Cells.Item(3,5).Value = 650

# Row 4 - BRVM - AUTRES SECTEURS
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 4
$ws.Cells.Item(4,4).Value = 2548.68
$ws.Cells.Item(4,5).Value = 640.71

# Row 5 - now NEI-CEDA CI (was SAFCA CI)
$ws.Cells.Item(5,1).Value = "NEI-CEDA CI"
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 4
$ws.Cells.Item(5,4).Value = 2375
$ws.Cells.Item(5,5).Value = 595

# Row 6 - now SETAO CI (was NEI-CEDA CI)
$ws.Cells.Item(6,1).Value = "SETAO CI"
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 4
$ws.Cells.Item(6,4).Value = 2310
$ws.Cells.Item(6,5).Value = 560

# Row 7 - now UNIWAX CI (was SETAO CI)
$ws.Cells.Item(7,1).Value = "UNIWAX CI"
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 4
$ws.Cells.Item(7,4).Value = 2310
$ws.Cells.Item(7,5).Value = 585

# Row 8 - now AIR LIQUIDE CI (was UNIWAX CI)
$ws.Cells.Item(8,1).Value = "AIR LIQUIDE CI"
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(8,3).Value = 4
$ws.Cells.Item(8,4).Value = 2130
$ws.Cells.Item(8,5).Value = 540

# Row 9 - now SAFCA CI (was AIR LIQUIDE CI)
$ws.Cells.Item(9,1).Value = "SAFCA CI"
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,3).Value = 2
$ws.Cells.Item(9,4).Value = 1905
$ws.Cells.Item(9,5).Value = 985

# Row 10 - BRVM - DISTRIBUTION
$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(10,3).Value = 4
$ws.Cells.Item(10,4).Value = 1457.14
$ws.Cells.Item(10,5).Value = 367.96

# Row 11 - BRVM - TRANSPORT
$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 4
$ws.Cells.Item(11,4).Value = 1402.54
$ws.Cells.Item(11,5).Value = 352.47

# Row 12 - BRVM - AGRICULTURE
$ws.Cells.Item(12,2).Value = 0
$ws.Cells.Item(12,3).Value = 4
$ws.Cells.Item(12,4).Value = 1278.76
$ws.Cells.Item(12,5).Value = 323.48

# Row 13 - BRVM - INDUSTRIE
$ws.Cells.Item(13,2).Value = 0
$ws.Cells.Item(13,3).Value = 4
$ws.Cells.Item(13,4).Value = 1049.14
$ws.Cells.Item(13,5).Value = 262.99

# Row 14 - BRVM - CONSOMMATION DE BASE
$ws.Cells.Item(14,2).Value = 0
$ws.Cells.Item(14,3).Value = 4
$ws.Cells.Item(14,4).Value = 865.73
$ws.Cells.Item(14,5).Value = 217.65

# Row 15 - BRVM-PRINCIPAL
$ws.Cells.Item(15,2).Value = 0
$ws.Cells.Item(15,3).Value = 4
$ws.Cells.Item(15,4).Value = 753.04
$ws.Cells.Item(15,5).Value = 189.11

# Row 16 - BRVM - INDUSTRIELS
$ws.Cells.Item(16,2).Value = 0
$ws.Cells.Item(16,3).Value = 4
$ws.Cells.Item(16,4).Value = 550.98
$ws.Cells.Item(16,5).Value = 137.32

# Row 17 - BRVM-PRESTIGE
$ws.Cells.Item(17,2).Value = 0
$ws.Cells.Item(17,3).Value = 4
$ws.Cells.Item(17,4).Value = 514.05
$ws.Cells.Item(17,5).Value = 127.91

# Row 18 - BRVM - FINANCES
$ws.Cells.Item(18,2).Value = 0
$ws.Cells.Item(18,3).Value = 4
$ws.Cells.Item(18,4).Value = 483.94
$ws.Cells.Item(18,5).Value = 120.64

# Row 19 - BRVM - SERVICES FINANCIERS
$ws.Cells.Item(19,2).Value = 0
$ws.Cells.Item(19,3).Value = 4
$ws.Cells.Item(19,4).Value = 475.6
$ws.Cells.Item(19,5).Value = 118.56

# Row 20 - BRVM - ENERGIE
$ws.Cells.Item(20,2).Value = 0
$ws.Cells.Item(20,3).Value = 4
$ws.Cells.Item(20,4).Value = 432.86
$ws.Cells.Item(20,5).Value = 109.47

# Row 21 - BRVM - CONSOMMATION DISCRETIONNAIRE
$ws.Cells.Item(21,2).Value = 0
$ws.Cells.Item(21,3).Value = 4
$ws.Cells.Item(21,4).Value = 417.45
$ws.Cells.Item(21,5).Value = 104.87

# Row 22 - BRVM - TELECOMMUNICATIONS
$ws.Cells.Item(22,2).Value = 0
$ws.Cells.Item(22,3).Value = 4
$ws.Cells.Item(22,4).Value = 367.85
$ws.Cells.Item(22,5).Value = 92.11

# Row 23 - SAFCA CI (SAFC)
$ws.Cells.Item(23,2).Value = 4
$ws.Cells.Item(23,4).Value = 28.27

# Row 24 - now VIVO ENERGY CI (SHEC) (was BANK OF AFRICA ML (BOAM))
$ws.Cells.Item(24,1).Value = "VIVO ENERGY CI (SHEC)"
$ws.Cells.Item(24,2).Value = 2
$ws.Cells.Item(24,4).Value = 8.44
$ws.Cells.Item(24,5).Value = 5.09

# Row 25 - now BANK OF AFRICA ML (BOAM) (was SAPH CI (SPHC))
$ws.Cells.Item(25,1).Value = "BANK OF AFRICA ML (BOAM)"
$ws.Cells.Item(25,4).Value = 7.37
$ws.Cells.Item(25,5).Value = 7.37

# Row 26 - now BERNABE CI (BNBC) (was SICABLE CI (CABC))
$ws.Cells.Item(26,1).Value = "BERNABE CI (BNBC)"
$ws.Cells.Item(26,2).Value = 2
$ws.Cells.Item(26,3).Value = 1
$ws.Cells.Item(26,4).Value = 5.54
$ws.Cells.Item(26,5).Value = 2.19
$ws.Cells.Item(26,7).Value = "👀 À surveiller"

# Row 27 - now SAPH CI (SPHC) (was BERNABE CI (BNBC))
$ws.Cells.Item(27,1).Value = "SAPH CI (SPHC)"
$ws.Cells.Item(27,3).Value = 0
$ws.Cells.Item(27,4).Value = 4.69
$ws.Cells.Item(27,5).Value = 4.69
$ws.Cells.Item(27,7).Value = "➖ Neutre"

# Row 28 - now SERVAIR ABIDJAN CI (ABJC) (was VIVO ENERGY CI (SHEC))
$ws.Cells.Item(28,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws.Cells.Item(28,4).Value = 4.07
$ws.Cells.Item(28,5).Value = 4.07

# Row 29 - now SICABLE CI (CABC) (was AIR LIQUIDE CI (SIVC))
$ws.Cells.Item(29,1).Value = "SICABLE CI (CABC)"
$ws.Cells.Item(29,3).Value = 0
$ws.Cells.Item(29,4).Value = 4
$ws.Cells.Item(29,5).Value = 4
$ws.Cells.Item(29,7).Value = "➖ Neutre"

# Row 30 - now AIR LIQUIDE CI (SIVC) (was SICOR CI (SICC))
$ws.Cells.Item(30,1).Value = "AIR LIQUIDE CI (SIVC)"
$ws.Cells.Item(30,3).Value = 1
$ws.Cells.Item(30,4).Value = 2.94
$ws.Cells.Item(30,5).Value = -1.82
$ws.Cells.Item(30,7).Value = "👀 À surveiller"

# Row 31 - now SICOR CI (SICC) (was BANK OF AFRICA NG (BOAN))
$ws.Cells.Item(31,1).Value = "SICOR CI (SICC)"
$ws.Cells.Item(31,4).Value = 2.8
$ws.Cells.Item(31,5).Value = 2.8

# Row 32 - now NEI-CEDA CI (NEIC) (was UNIWAX CI (UNXC))
$ws.Cells.Item(32,1).Value = "NEI-CEDA CI (NEIC)"
$ws.Cells.Item(32,4).Value = 1.69
$ws.Cells.Item(32,5).Value = 1.69

# Row 33 - now ECOBANK TRANS. INCORP. TG (ETIT) (was NEI-CEDA CI (NEIC))
$ws.Cells.Item(33,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws.Cells.Item(33,3).Value = 1
$ws.Cells.Item(33,4).Value = 0.42
$ws.Cells.Item(33,5).Value = 6.67
$ws.Cells.Item(33,7).Value = "👀 À surveiller"

# Row 34 - now TOTAL (was ECOBANK TRANS. INCORP. TG (ETIT))
$ws.Cells.Item(34,1).Value = "TOTAL"
$ws.Cells.Item(34,2).Value = 0
$ws.Cells.Item(34,3).Value = 4
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(34,7).Value = "➖ Neutre"

# Row 35 - now BANK OF AFRICA NG (BOAN) (was TOTAL)
$ws.Cells.Item(35,1).Value = "BANK OF AFRICA NG (BOAN)"
$ws.Cells.Item(35,2).Value = 1
$ws.Cells.Item(35,3).Value = 1
$ws.Cells.Item(35,4).Value = -0.47
$ws.Cells.Item(35,5).Value = -2.71
$ws.Cells.Item(35,7).Value = "👀 À surveiller"

# Rows 36-39 are unchanged.

# Row 40 - now FILTISAC CI (FTSC) (was ONATEL BF (ONTBF))
$ws.Cells.Item(40,1).Value = "FILTISAC CI (FTSC)"
$ws.Cells.Item(40,4).Value = -1.84
$ws.Cells.Item(40,5).Value = -1.84

# Row 41 - now ONATEL BF (ONTBF) (was BICI CI (BICC))
$ws.Cells.Item(41,1).Value = "ONATEL BF (ONTBF)"
$ws.Cells.Item(41,4).Value = -2.13
$ws.Cells.Item(41,5).Value = -2.13

# Row 42 - now ORAGROUP TOGO (ORGT) (was SETAO CI (STAC))
$ws.Cells.Item(42,1).Value = "ORAGROUP TOGO (ORGT)"
$ws.Cells.Item(42,2).Value = 1
$ws.Cells.Item(42,3).Value = 2
$ws.Cells.Item(42,4).Value = -2.68
$ws.Cells.Item(42,5).Value = 2.53
$ws.Cells.Item(42,7).Value = "👀 À surveiller"

# Row 43 - now BICI CI (BICC) (was ORAGROUP TOGO (ORGT))
$ws.Cells.Item(43,1).Value = "BICI CI (BICC)"
$ws.Cells.Item(43,3).Value = 1
$ws.Cells.Item(43,4).Value = -2.85
$ws.Cells.Item(43,5).Value = -2.85

# Row 44 - new: UNIWAX CI (UNXC)
$ws.Cells.Item(44,1).Value = "UNIWAX CI (UNXC)"
$ws.Cells.Item(44,2).Value = 1
$ws.Cells.Item(44,3).Value = 1
$ws.Cells.Item(44,4).Value = -3.42
$ws.Cells.Item(44,5).Value = -5.17
$ws.Cells.Item(44,6).Value = "🟡 Observer"
$ws.Cells.Item(44,7).Value = "👀 À surveiller"

# Row 45 - new: SETAO CI (STAC)
$ws.Cells.Item(45,1).Value = "SETAO CI (STAC)"
$ws.Cells.Item(45,2).Value = 0
$ws.Cells.Item(45,3).Value = 1
$ws.Cells.Item(45,4).Value = -3.45
$ws.Cells.Item(45,5).Value = -3.45
$ws.Cells.Item(45,6).Value = "🟡 Observer"
$ws.Cells.Item(45,7).Value = "➖ Neutre"

# Row 46 - new: SOCIETE GENERALE COTE D'IVOIRE (SGBC)
$ws.Cells.Item(46,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws.Cells.Item(46,2).Value = 0
$ws.Cells.Item(46,3).Value = 1
$ws.Cells.Item(46,4).Value = -3.71
$ws.Cells.Item(46,5).Value = -3.71
$ws.Cells.Item(46,6).Value = "🟡 Observer"
$ws.Cells.Item(46,7).Value = "➖ Neutre"

# Row 47 - new: CIE CI (CIEC) (previously row 44, now re-computed & moved)
$ws.Cells.Item(47,1).Value = "CIE CI (CIEC)"
$ws.Cells.Item(47,2).Value = 0
$ws.Cells.Item(47,3).Value = 3
$ws.Cells.Item(47,4).Value = -8.89
$ws.Cells.Item(47,5).Value = -2.27
$ws.Cells.Item(47,6).Value = "🔴 Vente"
$ws.Cells.Item(47,7).Value = "⚠️ Risque de décrochage"

# ---------------------------------------------------------------------
# Sheet 2: Top_YTD
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Cells.Item(2,2).Value = 6853884.24

$ws2.Cells.Item(3,2).Value = 314196.88

$ws2.Cells.Item(4,2).Value = 294902.1

$ws2.Cells.Item(5,2).Value = 231534.64

# Row 6 - now UNIWAX CI (was SETAO CI)
$ws2.Cells.Item(6,1).Value = "UNIWAX CI"
$ws2.Cells.Item(6,2).Value = 210558.05

# Row 7 - now SETAO CI (was UNIWAX CI)
$ws2.Cells.Item(7,1).Value = "SETAO CI"
$ws2.Cells.Item(7,2).Value = 210476.96

$ws2.Cells.Item(8,2).Value = 159860

# Row 9 - now BRVM - DISTRIBUTION (was SAFCA CI)
$ws2.Cells.Item(9,1).Value = "BRVM - DISTRIBUTION"
$ws2.Cells.Item(9,2).Value = 46364.02

# Row 10 - now BRVM - TRANSPORT (was BRVM - DISTRIBUTION)
$ws2.Cells.Item(10,1).Value = "BRVM - TRANSPORT"
$ws2.Cells.Item(10,2).Value = 41134.07

# Row 11 - now BRVM - AGRICULTURE (was BRVM - TRANSPORT)
$ws2.Cells.Item(11,1).Value = "BRVM - AGRICULTURE"
$ws2.Cells.Item(11,2).Value = 30921.16
